# Steuertarife.xlsx — "add income tax calculator"
#
# The underlying data change: on the "Staatssteuer" sheet, the
# "SteuerfussKanton" column (H) for rows 86-113 drops from a flat 100
# (i.e. 100%) down to 1 (i.e. 1%) for every one of those rows.
#
# The workbook was also left scrolled/selected further down by the
# author while they made the edit (topLeftCell A84->A95, selection
# I105->N90); we reproduce the selection change too.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Staatssteuer")

# H86:H113 — SteuerfussKanton 100 -> 1
$rng = $ws.Range("H86:H113")
$rng.Value = 1

# Reflect the author's final cursor position/selection on that sheet.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 95
$ws.Range("N90").Select()
